$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column M: "cpf_api" header + 10 CPF values -----------------------
$ws.Range("M1").Value = "cpf_api"
$ws.Range("M2").Value  = "'97093236014"
$ws.Range("M3").Value  = "'60094146012"
$ws.Range("M4").Value  = "'84809766080"
$ws.Range("M5").Value  = "'62648716050"
$ws.Range("M6").Value  = "'26276298085"
$ws.Range("M7").Value  = "'01317496094"
$ws.Range("M8").Value  = "'55856777050"
$ws.Range("M9").Value  = "'19626829001"
$ws.Range("M10").Value = "'24094592008"
$ws.Range("M11").Value = "'58063164083"

# --- Column M width ---------------------------------------------------------
$ws.Columns.Item(13).ColumnWidth = 12

# --- View / selection state --------------------------------------------------
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 3
$win.ScrollRow = 1
$excel.Goto($ws.Range("C1"), $true)

$win.Left = 615
$win.Top = 450
$win.Width = 18030
$win.Height = 9825

$ws.Range("M11").Select()
